# Append the new "2021/08/30" report row (row 74) to Sheet1, mirroring
# the existing row 73 exactly (same styles) but with the new date label
# and updated metrics, per DGS's 2021/08/30 report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone row 73's formatting/values down into row 74 so the new row
#    picks up the exact same cell styles (s="1" date-style text cell in
#    column A, s="2" numeric cells in B:E) without Excel inventing any
#    new style entries.
$ws.Range("A73:E73").Copy() | Out-Null
$ws.Range("A74:E74").PasteSpecial() | Out-Null

# 2) Column A stores these dates as literal text (not real dates), but
#    typing "2021/08/30" straight into a date-formatted cell would get
#    auto-parsed into a date serial. Route it through a scratch formula
#    cell instead: a formula's cached string result isn't subject to
#    that auto-detection, so copying *its value only* into A74 yields a
#    plain text cell - exactly like the rest of column A.
$scratch = $ws.Cells.Item(80, 1)
$scratch.Formula = "=""2021/08/30"""
$scratch.Copy() | Out-Null
$ws.Cells.Item(74, 1).PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null

# 3) Fill in the actual reported metrics for 2021/08/30.
$ws.Cells.Item(74, 2).Value = 297.7
$ws.Cells.Item(74, 3).Value = 303.3
$ws.Cells.Item(74, 4).Value = 0.98
$ws.Cells.Item(74, 5).Value = 0.99
